$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.918.47"
$ws.Range("D3").Value = "2.544.95"
$ws.Range("E3").Value = "  -5.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.93%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.26"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.81%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.72"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.76%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.116"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +7.58%  "
$ws.Range("D14").Value = "2.931.56"
$ws.Range("E14").Value = "  -5.35%  "
$ws.Range("D15").Value = "2.579.82"
$ws.Range("E15").Value = "  -4.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.879"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.20"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.36%  "
$ws.Range("D18").Value = "42.896.39"
$ws.Range("E18").Value = "  -5.58%  "
$ws.Range("D19").Value = "0.0₃0981"
$ws.Range("E19").Value = "  -3.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.67"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.54"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.84"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.72"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -9.47%  "
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.03"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -7.05%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.22"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.97%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.71"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.06"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.36"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.14"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -8.90%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.37"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -11.05%  "
$ws.Range("E36").Value = "  -6.19%  "
$ws.Range("E37").Value = "  -5.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.78"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.54%  "
$ws.Range("E39").Value = "  -4.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.98"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -10.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0310"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.85"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("D44").Value = "2.095.12"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.99"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +24.82%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -9.85%  "
$ws.Range("D49").Value = "2.788.57"
$ws.Range("E49").Value = "  -5.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.81"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.68"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.63%  "
